$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = 'Última actualización: 19:56:22'
$ws.Cells.Item(3, 1).Value = 'Total filas: 510'
$ws.Cells.Item(323, 1).Value = '14:46:52'
$ws.Cells.Item(323, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(323, 4).Value = 94
$ws.Cells.Item(324, 1).Value = '15:31:33'
$ws.Cells.Item(324, 3).Value = '215C_EL PATO'
$ws.Cells.Item(324, 4).Value = 49
$ws.Cells.Item(353, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(355, 3).Value = '215A_EL PATO'
$ws.Cells.Item(452, 1).Value = '18:37:25'
$ws.Cells.Item(452, 3).Value = '14_ABASTO'
$ws.Cells.Item(452, 4).Value = 44
$ws.Cells.Item(453, 1).Value = '17:41:19'
$ws.Cells.Item(453, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(453, 4).Value = 100
$ws.Cells.Item(454, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(477, 1).Value = '19:56:21'
$ws.Cells.Item(477, 2).Value = '19:57'
$ws.Cells.Item(477, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(477, 4).Value = 1
$ws.Cells.Item(478, 1).Value = '18:17:05'
$ws.Cells.Item(478, 2).Value = '19:59'
$ws.Cells.Item(478, 4).Value = 102
$ws.Cells.Item(479, 1).Value = '18:37:25'
$ws.Cells.Item(479, 2).Value = '20:00'
$ws.Cells.Item(479, 3).Value = '17_ROMERO'
$ws.Cells.Item(479, 4).Value = 83
$ws.Cells.Item(480, 1).Value = '19:56:21'
$ws.Cells.Item(480, 2).Value = '20:00'
$ws.Cells.Item(480, 3).Value = '14_ABASTO'
$ws.Cells.Item(480, 4).Value = 4
$ws.Cells.Item(481, 1).Value = '19:42:02'
$ws.Cells.Item(481, 2).Value = '20:01'
$ws.Cells.Item(481, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(481, 4).Value = 19
$ws.Cells.Item(482, 1).Value = '19:42:02'
$ws.Cells.Item(482, 2).Value = '20:09'
$ws.Cells.Item(482, 3).Value = '15_ABASTO'
$ws.Cells.Item(482, 4).Value = 27
$ws.Cells.Item(483, 2).Value = '20:10'
$ws.Cells.Item(483, 3).Value = '15_ABASTO'
$ws.Cells.Item(483, 4).Value = 79
$ws.Cells.Item(484, 1).Value = '19:56:21'
$ws.Cells.Item(484, 2).Value = '20:10'
$ws.Cells.Item(484, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(484, 4).Value = 14
$ws.Cells.Item(485, 1).Value = '18:17:05'
$ws.Cells.Item(485, 2).Value = '20:11'
$ws.Cells.Item(485, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(485, 4).Value = 114
$ws.Cells.Item(486, 1).Value = '18:51:07'
$ws.Cells.Item(486, 2).Value = '20:12'
$ws.Cells.Item(486, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(486, 4).Value = 81
$ws.Cells.Item(487, 1).Value = '19:42:02'
$ws.Cells.Item(487, 2).Value = '20:13'
$ws.Cells.Item(487, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(487, 4).Value = 31
$ws.Cells.Item(488, 1).Value = '18:58:44'
$ws.Cells.Item(488, 2).Value = '20:21'
$ws.Cells.Item(488, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(488, 4).Value = 83
$ws.Cells.Item(489, 1).Value = '19:42:02'
$ws.Cells.Item(489, 2).Value = '20:22'
$ws.Cells.Item(489, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(489, 4).Value = 40
$ws.Cells.Item(490, 1).Value = '18:37:25'
$ws.Cells.Item(490, 2).Value = '20:22'
$ws.Cells.Item(490, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(490, 4).Value = 105
$ws.Cells.Item(491, 1).Value = '18:37:25'
$ws.Cells.Item(491, 2).Value = '20:23'
$ws.Cells.Item(491, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(491, 4).Value = 106
$ws.Cells.Item(492, 1).Value = '19:56:21'
$ws.Cells.Item(492, 2).Value = '20:23'
$ws.Cells.Item(492, 3).Value = '215A_EL PATO'
$ws.Cells.Item(492, 4).Value = 27
$ws.Cells.Item(493, 1).Value = '18:37:25'
$ws.Cells.Item(493, 2).Value = '20:24'
$ws.Cells.Item(493, 3).Value = '215A_EL PATO'
$ws.Cells.Item(493, 4).Value = 107
$ws.Cells.Item(494, 1).Value = '18:51:07'
$ws.Cells.Item(494, 2).Value = '20:25'
$ws.Cells.Item(494, 3).Value = '215A_EL PATO'
$ws.Cells.Item(494, 4).Value = 94
$ws.Cells.Item(495, 1).Value = '19:42:02'
$ws.Cells.Item(495, 2).Value = '20:26'
$ws.Cells.Item(495, 4).Value = 44
$ws.Cells.Item(496, 1).Value = '18:51:07'
$ws.Cells.Item(496, 2).Value = '20:27'
$ws.Cells.Item(496, 3).Value = '14_ABASTO'
$ws.Cells.Item(496, 4).Value = 96
$ws.Cells.Item(497, 1).Value = '19:42:02'
$ws.Cells.Item(497, 2).Value = '20:31'
$ws.Cells.Item(497, 3).Value = '225_GOMEZ'
$ws.Cells.Item(497, 4).Value = 49
$ws.Cells.Item(498, 1).Value = '18:37:25'
$ws.Cells.Item(498, 2).Value = '20:32'
$ws.Cells.Item(498, 3).Value = '225_GOMEZ'
$ws.Cells.Item(498, 4).Value = 115
$ws.Cells.Item(499, 1).Value = '18:58:44'
$ws.Cells.Item(499, 2).Value = '20:35'
$ws.Cells.Item(499, 3).Value = '14_ABASTO'
$ws.Cells.Item(499, 4).Value = 97
$ws.Cells.Item(500, 2).Value = '20:39'
$ws.Cells.Item(500, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(500, 4).Value = 57
$ws.Cells.Item(501, 1).Value = '18:51:07'
$ws.Cells.Item(501, 2).Value = '20:46'
$ws.Cells.Item(501, 3).Value = '14X44_ABASTO'
$ws.Cells.Item(501, 4).Value = 115
$ws.Cells.Item(502, 1).Value = '18:58:44'
$ws.Cells.Item(502, 2).Value = '20:48'
$ws.Cells.Item(502, 3).Value = '14X44_ABASTO'
$ws.Cells.Item(502, 4).Value = 110
$ws.Cells.Item(503, 1).Value = '19:56:21'
$ws.Cells.Item(503, 2).Value = '20:52'
$ws.Cells.Item(503, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(503, 4).Value = 56
$ws.Cells.Item(504, 2).Value = '20:52'
$ws.Cells.Item(504, 3).Value = '15_ABASTO'
$ws.Cells.Item(504, 4).Value = 70
$ws.Cells.Item(505, 2).Value = '20:53'
$ws.Cells.Item(505, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(505, 4).Value = 71
$ws.Cells.Item(506, 1).Value = '18:58:44'
$ws.Cells.Item(506, 2).Value = '20:56'
$ws.Cells.Item(506, 4).Value = 118
$ws.Cells.Item(507, 2).Value = '20:57'
$ws.Cells.Item(507, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(507, 4).Value = 75
$ws.Cells.Item(508, 1).Value = '19:56:21'
$ws.Cells.Item(508, 2).Value = '21:00'
$ws.Cells.Item(508, 3).Value = '215B_EL PATO'
$ws.Cells.Item(508, 4).Value = 64
$ws.Cells.Item(509, 1).Value = '19:42:02'
$ws.Cells.Item(509, 2).Value = '21:01'
$ws.Cells.Item(509, 3).Value = '215B_EL PATO'
$ws.Cells.Item(509, 4).Value = 79
$ws.Cells.Item(509, 5).Value = 'LP1912'
$ws.Cells.Item(510, 1).Value = '19:42:02'
$ws.Cells.Item(510, 2).Value = '21:04'
$ws.Cells.Item(510, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(510, 4).Value = 82
$ws.Cells.Item(510, 5).Value = 'LP1912'
$ws.Cells.Item(511, 1).Value = '19:42:02'
$ws.Cells.Item(511, 2).Value = '21:21'
$ws.Cells.Item(511, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(511, 4).Value = 99
$ws.Cells.Item(511, 5).Value = 'LP1912'
$ws.Cells.Item(512, 1).Value = '19:42:02'
$ws.Cells.Item(512, 2).Value = '21:23'
$ws.Cells.Item(512, 3).Value = '10_OLMOS'
$ws.Cells.Item(512, 4).Value = 101
$ws.Cells.Item(512, 5).Value = 'LP1912'
$ws.Cells.Item(513, 1).Value = '19:42:02'
$ws.Cells.Item(513, 2).Value = '21:38'
$ws.Cells.Item(513, 3).Value = '14_ABASTO'
$ws.Cells.Item(513, 4).Value = 116
$ws.Cells.Item(513, 5).Value = 'LP1912'
$ws.Cells.Item(514, 1).Value = '19:42:02'
$ws.Cells.Item(514, 2).Value = '21:38'
$ws.Cells.Item(514, 3).Value = '17_ROMERO'
$ws.Cells.Item(514, 4).Value = 116
$ws.Cells.Item(514, 5).Value = 'LP1912'
$ws.Cells.Item(515, 1).Value = '19:56:21'
$ws.Cells.Item(515, 2).Value = '21:47'
$ws.Cells.Item(515, 3).Value = '215A_EL PATO'
$ws.Cells.Item(515, 4).Value = 111
$ws.Cells.Item(515, 5).Value = 'LP1912'

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = 'Última actualización: 19:56:22'
$ws.Cells.Item(3, 1).Value = 'Total filas: 55'
$ws.Cells.Item(55, 1).Value = '19:56:21'
$ws.Cells.Item(55, 2).Value = '20:23'
$ws.Cells.Item(55, 4).Value = 27
$ws.Cells.Item(56, 1).Value = '18:37:25'
$ws.Cells.Item(56, 2).Value = '20:24'
$ws.Cells.Item(56, 4).Value = 107
$ws.Cells.Item(57, 1).Value = '18:51:07'
$ws.Cells.Item(57, 2).Value = '20:25'
$ws.Cells.Item(57, 3).Value = '215A_EL PATO'
$ws.Cells.Item(57, 4).Value = 94
$ws.Cells.Item(58, 1).Value = '19:56:21'
$ws.Cells.Item(58, 2).Value = '21:00'
$ws.Cells.Item(58, 3).Value = '215B_EL PATO'
$ws.Cells.Item(58, 4).Value = 64
$ws.Cells.Item(58, 5).Value = 'LP1912'
$ws.Cells.Item(59, 1).Value = '19:42:02'
$ws.Cells.Item(59, 2).Value = '21:01'
$ws.Cells.Item(59, 3).Value = '215B_EL PATO'
$ws.Cells.Item(59, 4).Value = 79
$ws.Cells.Item(59, 5).Value = 'LP1912'
$ws.Cells.Item(60, 1).Value = '19:56:21'
$ws.Cells.Item(60, 2).Value = '21:47'
$ws.Cells.Item(60, 3).Value = '215A_EL PATO'
$ws.Cells.Item(60, 4).Value = 111
$ws.Cells.Item(60, 5).Value = 'LP1912'

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = 'Última actualización: 19:56:22'
$ws.Cells.Item(3, 1).Value = 'Total filas: 67'
$ws.Cells.Item(67, 1).Value = '19:56:21'
$ws.Cells.Item(67, 2).Value = '20:02'
$ws.Cells.Item(67, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(67, 4).Value = 6
$ws.Cells.Item(67, 5).Value = 'L6203'
$ws.Cells.Item(68, 1).Value = '19:56:21'
$ws.Cells.Item(68, 2).Value = '20:39'
$ws.Cells.Item(68, 4).Value = 43
$ws.Cells.Item(69, 1).Value = '18:51:07'
$ws.Cells.Item(69, 2).Value = '20:40'
$ws.Cells.Item(69, 4).Value = 109
$ws.Cells.Item(70, 2).Value = '20:41'
$ws.Cells.Item(70, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(70, 4).Value = 59
$ws.Cells.Item(70, 5).Value = 'L6173'
$ws.Cells.Item(71, 1).Value = '19:17:03'
$ws.Cells.Item(71, 2).Value = '20:43'
$ws.Cells.Item(71, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(71, 4).Value = 86
$ws.Cells.Item(71, 5).Value = 'L6173'
$ws.Cells.Item(72, 1).Value = '19:42:02'
$ws.Cells.Item(72, 2).Value = '21:29'
$ws.Cells.Item(72, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(72, 4).Value = 107
$ws.Cells.Item(72, 5).Value = 'L6203'
